# Auto-generated Excel COM-interop script to apply cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($range, [string]$value) {
    # Force text entry so numeric-looking strings (e.g. "1.00", "7.66")
    # are not silently coerced into numbers by Excel, while leaving the
    # cell style untouched (matches original formatting).
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2
Set-CellText $ws.Range("D2") "66.437.63"
Set-CellText $ws.Range("E2") "  -1.66%  "

# Row 3
Set-CellText $ws.Range("D3") "2.509.32"
Set-CellText $ws.Range("E3") "  -4.58%  "

# Row 4
Set-CellText $ws.Range("D4") "1.00"
Set-CellText $ws.Range("E4") "  -0.01%  "

# Row 5
Set-CellText $ws.Range("D5") "582.13"
Set-CellText $ws.Range("E5") "  -2.12%  "

# Row 6
Set-CellText $ws.Range("D6") "172.68"
Set-CellText $ws.Range("E6") "  +2.43%  "

# Row 7
Set-CellText $ws.Range("E7") "  +0.03%  "

# Row 8
Set-CellText $ws.Range("E8") "  -2.48%  "

# Row 9
Set-CellText $ws.Range("D9") "2.506.52"
Set-CellText $ws.Range("E9") "  -4.71%  "

# Row 10
Set-CellText $ws.Range("E10") "  -0.77%  "

# Row 11
Set-CellText $ws.Range("E11") "  -0.29%  "

# Row 12
Set-CellText $ws.Range("D12") "0.350"
Set-CellText $ws.Range("E12") "  -3.97%  "

# Row 13
Set-CellText $ws.Range("D13") "5.11"
Set-CellText $ws.Range("E13") "  -2.21%  "

# Row 14
Set-CellText $ws.Range("D14") "26.53"
Set-CellText $ws.Range("E14") "  -4.21%  "

# Row 15
Set-CellText $ws.Range("D15") "2.940.32"
Set-CellText $ws.Range("E15") "  -5.44%  "

# Row 16
Set-CellText $ws.Range("E16") "  -3.85%  "

# Row 17
Set-CellText $ws.Range("D17") "66.347.99"
Set-CellText $ws.Range("E17") "  -1.55%  "

# Row 18
Set-CellText $ws.Range("D18") "2.469.15"
Set-CellText $ws.Range("E18") "  -6.62%  "

# Row 19
Set-CellText $ws.Range("D19") "11.25"
Set-CellText $ws.Range("E19") "  -6.56%  "

# Row 20
Set-CellText $ws.Range("D20") "7.66"
Set-CellText $ws.Range("E20") "  -5.09%  "

# Row 21
Set-CellText $ws.Range("D21") "346.50"
Set-CellText $ws.Range("E21") "  -3.11%  "

# Row 22
Set-CellText $ws.Range("E22") "  -3.01%  "

# Row 23
Set-CellText $ws.Range("D23") "4.59"
Set-CellText $ws.Range("E23") "  -1.84%  "

# Row 24
Set-CellText $ws.Range("B24") "SuiNetwork"
Set-CellText $ws.Range("C24") "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-CellText $ws.Range("D24") "1.95"
Set-CellText $ws.Range("E24") "  +0.07%  "

# Row 25
Set-CellText $ws.Range("B25") "Dai"
Set-CellText $ws.Range("C25") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-CellText $ws.Range("D25") "1.00"
Set-CellText $ws.Range("E25") "  +0.21%  "

# Row 26
Set-CellText $ws.Range("D26") "69.31"
Set-CellText $ws.Range("E26") "  -0.81%  "

# Row 27
Set-CellText $ws.Range("D27") "9.85"
Set-CellText $ws.Range("E27") "  -4.64%  "

# Row 28
Set-CellText $ws.Range("D28") "1.00"
Set-CellText $ws.Range("E28") "  +0.00%  "

# Row 29
Set-CellText $ws.Range("D29") "2.630.06"
Set-CellText $ws.Range("E29") "  -4.67%  "

# Row 30
Set-CellText $ws.Range("D30") "0.0₃0973"
Set-CellText $ws.Range("E30") "  -3.53%  "

# Row 31
Set-CellText $ws.Range("D31") "530.19"
Set-CellText $ws.Range("E31") "  -3.17%  "

# Row 32
Set-CellText $ws.Range("D32") "8.08"
Set-CellText $ws.Range("E32") "  +1.73%  "

# Row 33
Set-CellText $ws.Range("E33") "  -2.45%  "

# Row 34
Set-CellText $ws.Range("E34") "  -3.55%  "

# Row 35
Set-CellText $ws.Range("E35") "  -3.54%  "

# Row 36
Set-CellText $ws.Range("D36") "0.999"
Set-CellText $ws.Range("E36") "  -0.12%  "

# Row 37
Set-CellText $ws.Range("D37") "157.73"
Set-CellText $ws.Range("E37") "  +0.53%  "

# Row 38
Set-CellText $ws.Range("E38") "  -3.93%  "

# Row 39
Set-CellText $ws.Range("D39") "18.53"

# Row 40
Set-CellText $ws.Range("E40") "  +0.26%  "

# Row 41
Set-CellText $ws.Range("E41") "  -3.63%  "

# Row 42
Set-CellText $ws.Range("D42") "1.79"
Set-CellText $ws.Range("E42") "  -1.88%  "

# Row 43
Set-CellText $ws.Range("E43") "  -2.95%  "

# Row 45
Set-CellText $ws.Range("D45") "2.47"
Set-CellText $ws.Range("E45") "  +1.90%  "

# Row 46
Set-CellText $ws.Range("D46") "147.17"
Set-CellText $ws.Range("E46") "  -3.75%  "

# Row 47
Set-CellText $ws.Range("D47") "0.556"
Set-CellText $ws.Range("E47") "  -4.13%  "

# Row 48
Set-CellText $ws.Range("D48") "3.67"
Set-CellText $ws.Range("E48") "  -3.44%  "

# Row 49
Set-CellText $ws.Range("D49") "1.73"
Set-CellText $ws.Range("E49") "  +1.72%  "

# Row 50
Set-CellText $ws.Range("D50") "0.0₆0269"
Set-CellText $ws.Range("E50") "  -9.55%  "

# Row 51
Set-CellText $ws.Range("E51") "  -2.68%  "
